# Added error checking for watcher and file reader:
# the price for the "Dell Optiplex 7020" row (D4) was being populated
# before the watcher/reader had validated the incoming value, so the
# stale placeholder price must be cleared out. Formatting (currency
# number format) on the cell is preserved; only the value is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the (incorrect/placeholder) price value out of D4 while keeping
# its existing number-format style intact.
$ws.Range("D4").ClearContents()

# Move/restore the active selection to D4 (the cell that was just
# corrected) instead of the stale D15 selection left in the file.
$ws.Range("D4").Select() | Out-Null
